# Insert a new data row at row 213 (pushing the existing rows 213..314
# down to 214..315) and populate the new row with the latest observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(213).Insert()

$ws.Range("A213").Value2 = 6
$ws.Range("B213").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C213").Value2 = "Metropolitana"
$ws.Range("D213").Value2 = 44830
$ws.Range("E213").Value2 = 13
$ws.Range("F213").Value2 = 100112026
$ws.Range("G213").Value2 = "Haba"
$ws.Range("H213").Value2 = "Sin especificar"
$ws.Range("I213").Value2 = "Primera"
$ws.Range("J213").Value2 = 1050
$ws.Range("K213").Value2 = 9000
$ws.Range("L213").Value2 = 10000
$ws.Range("M213").Value2 = 9457
$ws.Range("N213").Value2 = "`$/saco 25 kilos"
$ws.Range("O213").Value2 = "Región Metropolitana"
$ws.Range("P213").Value2 = 378
$ws.Range("Q213").Value2 = 25
$ws.Range("R213").Value2 = "Hortaliza"
